# Auto-generated PowerShell COM-interop script
# Implements the commit: adds a "movement_date" column (AF) to the
# "All Cards in Done" sheet, and refreshes Trello label "uses" counters
# embedded inside the S column's JSON-like text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All Cards in Done")

# ---------------------------------------------------------------------------
# 1) Refresh label usage counters embedded in the "labels" (S) column text.
#    These are plain textual substitutions inside the inline-string JSON
#    blobs -- no structural change, just updated counters.
# ---------------------------------------------------------------------------
$lastRow = 132
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 19)   # column S
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $newVal = $val
        $newVal = $newVal.Replace("'uses': 369}", "'uses': 376}")
        $newVal = $newVal.Replace("'uses': 321}", "'uses': 322}")
        $newVal = $newVal.Replace("'uses': 86}", "'uses': 87}")
        $newVal = $newVal.Replace("'uses': 233}", "'uses': 238}")
        $newVal = $newVal.Replace("'uses': 196}", "'uses': 200}")
        if ($newVal -ne $val) {
            $cell.Value2 = $newVal
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Add the new "movement_date" column (AF) with header + per-row values.
# ---------------------------------------------------------------------------
$ws.Range("AF1").Value2 = "movement_date"

$movementDates = @(
    @{Row=2; Val=45321.64441990741},
    @{Row=3; Val=45321.4846211574},
    @{Row=4; Val=45320.86833048611},
    @{Row=5; Val=45320.74606923611},
    @{Row=6; Val=45320.73730489583},
    @{Row=7; Val=45320.6677017824},
    @{Row=8; Val=45317.66264984954},
    @{Row=9; Val=45317.65143436343},
    @{Row=10; Val=45317.59274125},
    @{Row=11; Val=45315.95919740741},
    @{Row=12; Val=45314.85763059028},
    @{Row=13; Val=45314.83411855324},
    @{Row=14; Val=45314.67953149306},
    @{Row=15; Val=45313.97050351852},
    @{Row=16; Val=45310.75418130787},
    @{Row=17; Val=45310.75352620371},
    @{Row=18; Val=45310.67675300926},
    @{Row=19; Val=45310.58195571759},
    @{Row=20; Val=45310.57904549769},
    @{Row=21; Val=45310.5128396875},
    @{Row=22; Val=45309.88025594907},
    @{Row=23; Val=45309.87965506945},
    @{Row=24; Val=45309.84264208333},
    @{Row=25; Val=45309.77911354167},
    @{Row=26; Val=45309.6701758912},
    @{Row=27; Val=45309.65183231481},
    @{Row=28; Val=45309.65076390046},
    @{Row=29; Val=45309.64983065972},
    @{Row=30; Val=45309.64704370371},
    @{Row=31; Val=45309.64321545139},
    @{Row=32; Val=45309.58682651621},
    @{Row=33; Val=45308.73196929398},
    @{Row=34; Val=45307.82388806713},
    @{Row=35; Val=45306.71879814815},
    @{Row=36; Val=45306.68796200232},
    @{Row=37; Val=45303.50352185185},
    @{Row=38; Val=45302.96849951389},
    @{Row=39; Val=45302.96611523148},
    @{Row=40; Val=45302.96159962963},
    @{Row=41; Val=45302.95861336806},
    @{Row=42; Val=45302.9448496412},
    @{Row=43; Val=45302.91460706019},
    @{Row=44; Val=45302.89107320602},
    @{Row=45; Val=45302.88931747685},
    @{Row=46; Val=45302.84969537037},
    @{Row=47; Val=45302.84855047454},
    @{Row=48; Val=45302.82785723379},
    @{Row=49; Val=45302.82566556713},
    @{Row=50; Val=45302.8169397801},
    @{Row=51; Val=45302.7340796412},
    @{Row=52; Val=45301.95290150463},
    @{Row=53; Val=45301.95218846065},
    @{Row=54; Val=45301.92891680555},
    @{Row=55; Val=45301.91510361111},
    @{Row=56; Val=45302.74915681713},
    @{Row=57; Val=45301.91164971065},
    @{Row=58; Val=45301.8640847338},
    @{Row=59; Val=45301.72878667824},
    @{Row=60; Val=45301.72718171297},
    @{Row=61; Val=45301.71266525463},
    @{Row=62; Val=45301.70071818287},
    @{Row=63; Val=45301.6853140625},
    @{Row=64; Val=45301.60337752315},
    @{Row=65; Val=45300.98007916666},
    @{Row=66; Val=45296.79376248843},
    @{Row=67; Val=45296.57727751158},
    @{Row=68; Val=45294.94193791666},
    @{Row=69; Val=45294.94122461806},
    @{Row=70; Val=45294.86376837963},
    @{Row=71; Val=45294.85913269676},
    @{Row=72; Val=45294.81978753473},
    @{Row=73; Val=45294.81969760417},
    @{Row=74; Val=45294.81958255787},
    @{Row=75; Val=45288.72397538194},
    @{Row=76; Val=45282.58166978009},
    @{Row=77; Val=45282.57716027777},
    @{Row=78; Val=45281.94473803241},
    @{Row=79; Val=45281.86092715277},
    @{Row=80; Val=45281.77926972222},
    @{Row=81; Val=45281.77924885417},
    @{Row=82; Val=45281.64254180556},
    @{Row=83; Val=45281.43172827546},
    @{Row=84; Val=45281.4215430324},
    @{Row=85; Val=45281.41371859954},
    @{Row=86; Val=45280.7736556713},
    @{Row=87; Val=45280.77282775463},
    @{Row=88; Val=45280.72932037037},
    @{Row=89; Val=45280.7215578125},
    @{Row=90; Val=45280.69797907407},
    @{Row=91; Val=45279.9587033912},
    @{Row=92; Val=45279.84811211805},
    @{Row=93; Val=45278.91231490741},
    @{Row=94; Val=45278.91229980324},
    @{Row=95; Val=45278.88792607639},
    @{Row=96; Val=45274.67148517361},
    @{Row=97; Val=45273.83550787037},
    @{Row=98; Val=45273.50947658565},
    @{Row=99; Val=45272.8678574537},
    @{Row=100; Val=45272.66151157407},
    @{Row=101; Val=45272.58665163194},
    @{Row=102; Val=45271.7538206713},
    @{Row=103; Val=45271.72398444444},
    @{Row=104; Val=45271.72200959491},
    @{Row=105; Val=45271.71038612269},
    @{Row=106; Val=45271.71029019676},
    @{Row=107; Val=45271.71018223379},
    @{Row=108; Val=45271.67947991898},
    @{Row=109; Val=45271.61266383102},
    @{Row=110; Val=45267.98361020834},
    @{Row=111; Val=45267.94140542824},
    @{Row=112; Val=45267.91566591436},
    @{Row=113; Val=45267.91441982639},
    @{Row=114; Val=45267.91422902778},
    @{Row=115; Val=45267.91249532408},
    @{Row=116; Val=45267.91070800926},
    @{Row=117; Val=45267.89615033565},
    @{Row=118; Val=45267.88008541667},
    @{Row=119; Val=45267.87417387732},
    @{Row=120; Val=45267.87207349537},
    @{Row=121; Val=45267.75843868055},
    @{Row=122; Val=45266.91523304398},
    @{Row=123; Val=45266.91514569445},
    @{Row=124; Val=45265.59286556713},
    @{Row=125; Val=45264.91138417824},
    @{Row=126; Val=45264.90938443287},
    @{Row=127; Val=45264.90421420139},
    @{Row=128; Val=45264.90345635416},
    @{Row=129; Val=45264.88809135417},
    @{Row=130; Val=45264.68942390046},
    @{Row=131; Val=45264.68928572917},
    @{Row=132; Val=45261.62664386574}
)

# Seed the very first data cell and establish the clean date/time style
# (lowercase format registers numFmtId 164, uppercase format then mutates
# that same style in-place to numFmtId 165 -- matching the authored file).
$firstRow = $movementDates[0].Row
$firstVal = $movementDates[0].Val
$firstCell = $ws.Cells.Item($firstRow, 32)   # column AF
$firstCell.Value2 = $firstVal
$firstCell.NumberFormat = "yyyy-mm-dd h:mm:ss"
$firstCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Remaining rows reuse the already-registered uppercase format so no extra
# style entries are created.
for ($i = 1; $i -lt $movementDates.Length; $i++) {
    $entry = $movementDates[$i]
    $cell = $ws.Cells.Item($entry.Row, 32)
    $cell.Value2 = $entry.Val
    $cell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

Write-Host "movement_date column added; label usage counters refreshed."
